$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 42 ("rock_cycle_title" / "THE ROCK CYCLE"),
# pushing existing rows 42+ down by one.
$ws.Rows.Item(42).Insert()

$ws.Range("B42").Value = "THE ROCK CYCLE"
$ws.Range("A42").Value = "rock_cycle_title"

# Update the sheet view to match the new selection/scroll position.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("A42").Select()
